$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row -> [old F value, new F value] (for validation/reference only; we just set new values)
$fUpdates = @{
    2  = 3375
    4  = 2453
    6  = 341
    7  = 1393
    8  = 1099
    10 = 518
    11 = 1164
    15 = 8625
    17 = 2486
    18 = 256
    19 = 251
    22 = 588
    24 = 1152
    26 = 2022
    27 = 2065
    28 = 62
    29 = 1759
    30 = 238
    31 = 1917
    33 = 32
    34 = 47
    38 = 303
    39 = 59
    40 = 245
    41 = 423
    42 = 495
    44 = 262
}

# sheet "全部类型" has two extra rows inserted before row 42/44 relative to "展览" (rows 8 and 18),
# shifting subsequent row numbers by +2 for rows after 18.
$fUpdatesAllTypes = @{
    2  = 3375
    4  = 2453
    6  = 341
    7  = 1393
    9  = 1099
    11 = 518
    12 = 1164
    15 = 8625
    17 = 2486
    19 = 256
    20 = 251
    23 = 588
    25 = 1152
    27 = 2022
    28 = 2065
    29 = 1759
    30 = 238
    31 = 1917
    33 = 32
    34 = 47
    38 = 303
    39 = 59
    40 = 245
    41 = 423
    46 = 495
    49 = 262
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    if ($name -eq "展览") {
        $updates = $fUpdates
    } else {
        $updates = $fUpdatesAllTypes
    }

    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }

    # G5: 218 (number) -> "已售罄" (text)
    $ws.Range("G5").Value = "已售罄"
}
